$wb = $excel.ActiveWorkbook
$src = $wb.Worksheets.Item("Data Harian - Table")

# Copy the daily data table (header row + 29 data rows) and paste it into a
# brand-new sheet placed right after the source sheet.
$src.Range("A9:K37").Copy()

$new = $wb.Worksheets.Add($null, $src)
$new.Name = "Sheet1"
$new.Range("A1").PasteSpecial()

$new.Range("A1:K29").Select()
$new.Activate()
